# Apply the "add 2022-Q4 data" edit:
#  1. Insert a brand-new worksheet "2022-Q4" right after "总计", before "2022-Q3".
#  2. Populate it with the fund-holdings table for the new quarter.
#  3. Insert a new row 2 in "总计" and fill it with the 2022-Q4 summary figures,
#     re-numbering the sequential index column (A) for every data row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet, positioned right after "总计".
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($q3Sheet, $totalSheet)
$newSheet.Name = "2022-Q4"

# Header row (row 1) -- same column headers as every other quarterly sheet.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $newSheet.Cells.Item(1, $col).Value = $headers[$col - 2]
}

# Match the bold/bordered header style and the index-column style used by
# every other quarterly sheet: bold font, thin box border, center/top aligned.
$headerRange = $newSheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$indexRange = $newSheet.Range("A2:A10")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1

# Fund-holdings data for 2022-Q4 (9 funds), in the same order as the source diff.
# Columns: code, name, scale, stockPosition, positionRatio, marketValue, positionRank
$q4Funds = @(
    @("005585", "银河文体娱乐主题灵活配置混合A", "3.15", "88.90", "5.40", "0.1701", 4),
    @("013890", "国泰睿毅三年持有期混合A",       "4.82", "90.04", "3.37", "0.1624", 10),
    @("008602", "方正富邦新兴成长混合A",         "1.26", "87.26", "3.65", "0.0460", 6),
    @("005075", "富国研究量化精选混合",           "2.48", "90.71", "1.54", "0.0382", 6),
    @("003397", "银华体育文化灵活配置混合",       "0.53", "87.74", "4.53", "0.0240", 4),
    @("015667", "银河文体娱乐主题灵活配置混合C", "0.38", "88.90", "5.40", "0.0205", 4),
    @("013891", "国泰睿毅三年持有期混合C",       "0.45", "90.04", "3.37", "0.0152", 10),
    @("004250", "银河量化优选混合",               "0.30", "62.88", "1.37", "0.0041", 7),
    @("008603", "方正富邦新兴成长混合C",         "0.03", "87.26", "3.65", "0.0011", 6)
)

$rowIdx = 2
foreach ($fund in $q4Funds) {
    $newSheet.Cells.Item($rowIdx, 1).Value = $rowIdx - 2
    # Fund code looks numeric ("005585") -- force text so leading zeros survive.
    $newSheet.Cells.Item($rowIdx, 2).Value = "'" + $fund[0]
    # Fund name is never numeric-looking -- plain text assignment is safe.
    $newSheet.Cells.Item($rowIdx, 3).Value = $fund[1]
    # Scale / position / ratio / market value are numeric-looking text in the
    # source data (kept as text so trailing zeros like "5.40" are preserved).
    $newSheet.Cells.Item($rowIdx, 4).Value = "'" + $fund[2]
    $newSheet.Cells.Item($rowIdx, 5).Value = "'" + $fund[3]
    $newSheet.Cells.Item($rowIdx, 6).Value = "'" + $fund[4]
    $newSheet.Cells.Item($rowIdx, 7).Value = "'" + $fund[5]
    $newSheet.Cells.Item($rowIdx, 8).Value = $fund[6]
    # The leading-apostrophe text entries above pick up an automatic "quote
    # prefix" style; the source data rows carry no special formatting at all
    # (only the header row and column A do), so strip that back off again.
    $newSheet.Range("B" + $rowIdx + ":G" + $rowIdx).ClearFormats()
    $rowIdx++
}

# ---------------------------------------------------------------------------
# 2) Insert the new summary row into "总计" and shift the rest down.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# A fresh inserted row inherits the header's bold/centered formatting; the
# data rows in this table are unstyled, so strip it back off B2:D2.
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 9
$totalSheet.Cells.Item(2, 4).Value = 0.48

# The blank inserted row lost the index column's cell style -- copy it back
# from the row below (which still carries the original formatting) before
# writing values, so A2 matches the style used by every other A-column cell.
$totalSheet.Cells.Item(3, 1).Copy()
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122)

# Re-number the sequential index column (A2:A7 = 0..5) now that the table has
# grown by one row -- matches the freshly regenerated index in the source file.
for ($r = 2; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
